$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("You will be experimenting ...") lost its "_GoBack"
#    bookmark, and the two runs that used to straddle <w:bookmarkEnd/>
#    ("... Labs " and "are marked as ...") become a single merged run.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2Range = $p2.Range

$p2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="1CB0EBD0" w14:textId="5C7C2067" w:rsidR="001742A5" w:rsidRDefault="001742A5" w:rsidP="001742A5">' +
  '<w:r><w:t>You will be experimenting with various aspects of WICED Wi</w:t></w:r>' +
  '<w:r w:rsidR="008D3787"><w:t>-</w:t></w:r>' +
  '<w:r><w:t>Fi</w:t></w:r>' +
  '<w:r w:rsidR="008D3787"><w:t>&#174;</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> by completing the exercises below. Labs are marked as &#8220;Basic&#8221; and &#8220;Advanced&#8221;. You should make </w:t></w:r>' +
  '<w:r w:rsidR="00297C27"><w:t>sure you complete the b</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">asic exercises </w:t></w:r>' +
  '<w:r w:rsidR="00C13C5B"><w:t xml:space="preserve">first </w:t></w:r>' +
  '<w:r w:rsidR="00423E16"><w:t>and then work</w:t></w:r>' +
  '<w:r w:rsidR="00297C27"><w:t xml:space="preserve"> on the a</w:t></w:r>' +
  '<w:r><w:t>dvanced exercises as time allows.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$p2Range.InsertXML($p2Xml)

# ---------------------------------------------------------------------------
# 2) The signoff table's first header cell ("Initials") now shows a Wingdings
#    checkmark symbol instead of the word, and the "_GoBack" bookmark (which
#    used to sit in the paragraph above) now wraps that (now empty) spot.
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$cell = $tbl.Cell(1, 1)
$cellPara = $cell.Range.Paragraphs(1)
$cellParaRange = $cellPara.Range

$cellXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="6DFF5C89" w14:textId="47FF3C96" w:rsidR="006D3E13" w:rsidRPr="00264538" w:rsidRDefault="006D3E13" w:rsidP="00253927">' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' +
  '<w:rPr><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="18"/></w:rPr>' +
  '<w:sym w:font="Wingdings" w:char="F0FC"/></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# Remove the old bookmark first so the new one (inserted via XML below) is
# re-numbered back down to id 0 instead of being given the next free id.
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

$cellParaRange.InsertXML($cellXml)

Write-Host "Edit applied"
